$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I18").Value = 310.33334
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 310.33334
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -26.33334000000002
$ws.Range("N18").Value = -1068
$ws.Range("H64").Value = 3066.6667
$ws.Range("I64").Value = 3100
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3100
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2852
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3066.6667
$ws.Range("I67").Value = 3100
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3100
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2242
$ws.Range("N67").Value = -4716
$ws.Range("H129").Value = 1081.2273
$ws.Range("I129").Value = 565.38464
$ws.Range("J129").Value = 1297.5483
$ws.Range("K129").Value = 1696.15392
$ws.Range("L129").Value = 3892.6449
$ws.Range("M129").Value = 3303.84608
$ws.Range("N129").Value = -13892.6449
$ws.Range("H137").Value = 1327.05
$ws.Range("I137").Value = 903.86664
$ws.Range("K137").Value = 2711.59992
$ws.Range("M137").Value = -161.5999199999997
$ws.Range("H138").Value = 2447.6829
$ws.Range("I138").Value = 1811.0769
$ws.Range("J138").Value = 3551.1333
$ws.Range("K138").Value = 5433.2307
$ws.Range("L138").Value = 10653.3999
$ws.Range("M138").Value = -293.2307000000001
$ws.Range("N138").Value = -20933.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 466.66666
$ws.Range("I5").Value = 466.66666
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 466.66666
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -354.66666
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 323096.6
$ws.Range("I32").Value = 374103.66
$ws.Range("J32").Value = 13130.538
$ws.Range("K32").Value = 374103.66
$ws.Range("L32").Value = 13130.538
$ws.Range("M32").Value = -373816.66
$ws.Range("N32").Value = -13704.538
$ws.Range("H41").Value = 1939.5
$ws.Range("I41").Value = 1939.5
$ws.Range("K41").Value = 1939.5
$ws.Range("M41").Value = -1525.5
$ws.Range("H43").Value = 10447
$ws.Range("J43").Value = 13000
$ws.Range("L43").Value = 13000
$ws.Range("N43").Value = -13626

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 466.66666
$ws.Range("I4").Value = 466.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 466.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -351.66666
$ws.Range("N4").ClearContents()
$ws.Range("H64").Value = 604.3889
$ws.Range("I64").Value = 570.125
$ws.Range("J64").Value = 631.8
$ws.Range("K64").Value = 570.125
$ws.Range("L64").Value = 631.8
$ws.Range("M64").Value = -345.125
$ws.Range("N64").Value = -1081.8
$ws.Range("H67").Value = 604.3889
$ws.Range("I67").Value = 570.125
$ws.Range("J67").Value = 631.8
$ws.Range("K67").Value = 570.125
$ws.Range("L67").Value = 631.8
$ws.Range("M67").Value = 209.875
$ws.Range("N67").Value = -2191.8
$ws.Range("H130").Value = 95800
$ws.Range("J130").Value = 95800
$ws.Range("L130").Value = 95800
$ws.Range("N130").Value = -105840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6029.2646
$ws.Range("I31").Value = 1310.4375
$ws.Range("J31").Value = 10223.777
$ws.Range("K31").Value = 1310.4375
$ws.Range("L31").Value = 10223.777
$ws.Range("M31").Value = -1015.4375
$ws.Range("N31").Value = -10813.777
$ws.Range("H34").Value = 6029.2646
$ws.Range("I34").Value = 1310.4375
$ws.Range("J34").Value = 10223.777
$ws.Range("K34").Value = 1310.4375
$ws.Range("L34").Value = 10223.777
$ws.Range("M34").Value = -1108.4375
$ws.Range("N34").Value = -10627.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1178.2609
$ws.Range("H68").Value = 1245.1818
$ws.Range("J68").Value = 1469.6143
$ws.Range("L68").Value = 4408.8429
$ws.Range("N68").Value = -6030.8429
$ws.Range("H71").Value = 1245.1818
$ws.Range("J71").Value = 1469.6143
$ws.Range("L71").Value = 13226.5287
$ws.Range("N71").Value = -21338.5287
$ws.Range("H97").Value = 1934.6666
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 2602
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 7806
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -8798
$ws.Range("H107").Value = 1934.9048
$ws.Range("J107").Value = 2972.8845
$ws.Range("L107").Value = 8918.6535
$ws.Range("N107").Value = -12758.6535
$ws.Range("H113").Value = 812.46344
$ws.Range("I113").Value = 518.2083
$ws.Range("J113").Value = 1227.8823
$ws.Range("K113").Value = 1554.6249
$ws.Range("L113").Value = 3683.6469
$ws.Range("M113").Value = 615.3751
$ws.Range("N113").Value = -8023.6469
$ws.Range("H131").Value = 1160.2
$ws.Range("I131").Value = 1001
$ws.Range("J131").Value = 1266.3334
$ws.Range("K131").Value = 3003
$ws.Range("L131").Value = 3799.0002
$ws.Range("M131").Value = 2037
$ws.Range("N131").Value = -13879.0002
$ws.Range("H136").Value = 3275.0715
$ws.Range("I136").Value = 3090.889
$ws.Range("J136").Value = 3606.6
$ws.Range("K136").Value = 9272.667000000001
$ws.Range("L136").Value = 10819.8
$ws.Range("M136").Value = -4172.667000000001
$ws.Range("N136").Value = -21019.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2501.5557
$ws.Range("I126").Value = 2457
$ws.Range("K126").Value = 7371
$ws.Range("M126").Value = -4901
$ws.Range("H132").Value = 2080.5908
$ws.Range("I132").Value = 1397.4
$ws.Range("K132").Value = 4192.200000000001
$ws.Range("M132").Value = -1662.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 50005000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 50005000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 50005000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -50005590
$ws.Range("H61").Value = 5900.5
$ws.Range("I61").Value = 5860.8
$ws.Range("J61").Value = 5966.6665
$ws.Range("K61").Value = 5860.8
$ws.Range("L61").Value = 5966.6665
$ws.Range("M61").Value = -5658.8
$ws.Range("N61").Value = -6370.6665
$ws.Range("H113").Value = 5900.5
$ws.Range("I113").Value = 5860.8
$ws.Range("J113").Value = 5966.6665
$ws.Range("K113").Value = 5860.8
$ws.Range("L113").Value = 5966.6665
$ws.Range("M113").Value = -3690.8
$ws.Range("N113").Value = -10306.6665
$ws.Range("H132").Value = 3643.0535
$ws.Range("I132").Value = 3595.1292
$ws.Range("J132").Value = 3702.48
$ws.Range("K132").Value = 10785.3876
$ws.Range("L132").Value = 11107.44
$ws.Range("M132").Value = -8255.3876
$ws.Range("N132").Value = -16167.44

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 121821.18
$ws.Range("J29").Value = 215005.5
$ws.Range("L29").Value = 215005.5
$ws.Range("N29").Value = -215585.5
$ws.Range("H46").Value = 55214.5
$ws.Range("J46").Value = 55214.5
$ws.Range("L46").Value = 55214.5
$ws.Range("N46").Value = -55676.5
$ws.Range("H74").Value = 5000
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 5000
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 15000
$ws.Range("N77").Value = -24360
$ws.Range("H122").Value = 1857.7222
$ws.Range("I122").Value = 1510.3077
$ws.Range("J122").Value = 2761
$ws.Range("K122").Value = 4530.9231
$ws.Range("L122").Value = 8283
$ws.Range("M122").Value = -2080.9231
$ws.Range("N122").Value = -13183
$ws.Range("H134").Value = 55214.5
$ws.Range("J134").Value = 55214.5
$ws.Range("L134").Value = 165643.5
$ws.Range("N134").Value = -170713.5
